# Update "想去人数" (F) counts (and one status label in G) across all four
# sheets of the 上海-漫展信息 workbook, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 651
$ws.Range("F5").Value = 2910
$ws.Range("F6").Value = 19
$ws.Range("F7").Value = 236
$ws.Range("F10").Value = 6777
$ws.Range("F13").Value = 326
$ws.Range("F15").Value = 1473
$ws.Range("F17").Value = 1099
$ws.Range("F18").Value = 2202
$ws.Range("F19").Value = 1445
$ws.Range("F20").Value = 645
$ws.Range("F22").Value = 1087
$ws.Range("F23").Value = 85
$ws.Range("F24").Value = 157
$ws.Range("F26").Value = 1648
$ws.Range("F27").Value = 1633
$ws.Range("F31").Value = 1648
$ws.Range("F32").Value = 1180
$ws.Range("F35").Value = 17
$ws.Range("F37").Value = 390
$ws.Range("F38").Value = 2416
$ws.Range("F39").Value = 2680
$ws.Range("F41").Value = 176
$ws.Range("F43").Value = 12
$ws.Range("F44").Value = 16
$ws.Range("F45").Value = 308
$ws.Range("F46").Value = 117
$ws.Range("F47").Value = 156
$ws.Range("F48").Value = 135

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 6
$ws.Range("F7").Value = 143
$ws.Range("F15").Value = 53
$ws.Range("F17").Value = 159
$ws.Range("F19").Value = 42
$ws.Range("F20").Value = 16
$ws.Range("F23").Value = 451
$ws.Range("F24").Value = 38
$ws.Range("F31").Value = 8
$ws.Range("F37").Value = 37

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 534
$ws.Range("F6").Value = 1723
$ws.Range("F7").Value = 1650
$ws.Range("F8").Value = 1840
$ws.Range("F9").Value = 2688
$ws.Range("F10").Value = 979
$ws.Range("F11").Value = 882
$ws.Range("F13").Value = 227
$ws.Range("F14").Value = 1265
$ws.Range("F15").Value = 7026
$ws.Range("G15").Value = "已售罄"

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 534
$ws.Range("F4").Value = 1723
$ws.Range("F6").Value = 651
$ws.Range("F7").Value = 2910
$ws.Range("F8").Value = 236
$ws.Range("F9").Value = 1650
$ws.Range("F11").Value = 2688
$ws.Range("F12").Value = 6777
$ws.Range("F13").Value = 979
$ws.Range("F14").Value = 882
$ws.Range("F16").Value = 326
$ws.Range("F17").Value = 143
$ws.Range("F18").Value = 227
$ws.Range("F19").Value = 1099
$ws.Range("F20").Value = 2202
$ws.Range("F21").Value = 1445
$ws.Range("F22").Value = 645
$ws.Range("F24").Value = 1087
$ws.Range("F25").Value = 85
$ws.Range("F28").Value = 1648
$ws.Range("F29").Value = 159
$ws.Range("F31").Value = 1648
$ws.Range("F32").Value = 1180
$ws.Range("F35").Value = 451
$ws.Range("F36").Value = 390
$ws.Range("F37").Value = 38
$ws.Range("F39").Value = 2416
$ws.Range("F40").Value = 2680
$ws.Range("F42").Value = 176
$ws.Range("F43").Value = 308
$ws.Range("F44").Value = 117
$ws.Range("F45").Value = 156
